$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from the last existing data row (A33) down to
# the new rows so the new date cells keep the same numFmt (style index 2).
$ws.Range("A33").Copy()
$ws.Range("A34:A36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 34: 2010-03-23 (serial 40260), 2.5 hours, "Group Meeting"
$ws.Range("A34").Value = 40260
$ws.Range("B34").Value = 2.5
$ws.Range("C34").Value = "Group Meeting"

# Row 35: 2010-03-23 (serial 40260), 1 hour, "Weekly Meeting"
$ws.Range("A35").Value = 40260
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = "Weekly Meeting"

# Row 36: 2010-03-26 (serial 40263), 0.5 hours, "Skype Meeting"
$ws.Range("A36").Value = 40263
$ws.Range("B36").Value = 0.5
$ws.Range("C36").Value = "Skype Meeting"

$ws.Range("A37").Select()
